# Applies the two content changes from the diff:
#  1. The paragraph that holds the two inline/anchored figures (right
#     before the "<流れ図>" paragraph) has a stray empty <w:pPr><w:rPr>
#     <w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr> that should simply
#     disappear (the runs/drawings themselves are untouched).
#  2. The paragraph between "<流れ図>" and "<実行結果>" only contains a
#     leftover "_GoBack" bookmark plus the same stray <w:pPr><w:rPr>
#     <w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>; it should become a
#     totally empty paragraph (<w:p/>), i.e. the bookmark is deleted and
#     the paragraph-mark formatting is cleared.

$d = $word.ActiveDocument

# --- Paragraph 2: drop the leftover "_GoBack" bookmark -------------------
# Word keeps this one out of the regular Bookmarks collection/count (it's
# the special "last edit" bookmark) but it is still addressable by name.
$goBack = $null
try {
    $goBack = $d.Bookmarks("_GoBack")
} catch {
    $goBack = $null
}
if ($goBack -ne $null) {
    $goBack.Delete()
}

# --- Locate the two target paragraphs by their (now unique) content ------
$figurePara = $null
$blankPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text

    if ($figurePara -eq $null -and $para.Range.InlineShapes.Count -gt 0 -and $text -eq "") {
        $figurePara = $para
    } elseif ($blankPara -eq $null -and $text -eq "" -and $para.Range.InlineShapes.Count -eq 0 -and $para.Range.ShapeRange.Count -eq 0) {
        # the blank paragraph right after "<流れ図>" (the bookmark used to
        # live here); use the first still-unprocessed one we meet after
        # the figures paragraph has already been found.
        if ($figurePara -ne $null) {
            $blankPara = $para
        }
    }
}

# --- Paragraph 1: clear the stray paragraph-mark formatting --------------
# Re-applying the paragraph's own style forces Word to rebuild <w:pPr>
# from scratch, which drops the now-meaningless <w:rFonts w:hint="eastAsia"/>
# (the run-level formatting on the figures themselves is untouched).
if ($figurePara -ne $null) {
    $figurePara.Style = $figurePara.Style
}

# --- Paragraph 2: clear its paragraph-mark formatting too -----------------
if ($blankPara -ne $null) {
    $blankPara.Style = $blankPara.Style
}
